# Update the "userData" worksheet test data:
# - Row 3 ("invalidUserEmail"): remove the plot/address/contact columns (B:I),
#   keeping only the test case name (A) and the email (J) used by the test.
# - Row 4: turn it into a new "updateUser" test case row that only carries
#   the plot number (B), clearing the remaining address/contact columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:I3").ClearContents()

$ws.Range("C4:J4").ClearContents()
$ws.Range("A4").Value = "updateUser"

# Move the active selection, matching the saved workbook view state.
$ws.Range("D12").Select() | Out-Null
